$d = $word.ActiveDocument

# The two <id> elements were split across three runs each:
#   <id>            (Courier New, color 7f6000, sz 18)
#   p088v_aN         (plain, color 000000)
#   </id>           (Courier New, color 7f6000, sz 18)
# Collapse each trio down to a single run "<id>p088v_N</id>" that keeps the
# Courier-New/7f6000/18pt formatting of the surrounding tag runs.

$d.Content.Find.Execute("<id>p088v_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p088v_1</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p088v_a2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p088v_2</id>", 2) | Out-Null
